# Add the new "2022-Q4" quarterly sheet (copied/formatted like "2022-Q3"),
# positioned right after "总计", and refresh the "总计" summary sheet so its
# "2022-Q4" row is included (existing rows shift down by one).

$wb = $excel.ActiveWorkbook

$totalWs = $wb.Worksheets.Item(1)      # 总计
$q3Ws    = $wb.Worksheets.Item(2)      # 2022-Q3 (template for formatting)

# --- 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (keeps styles) ---
$q3Ws.Copy($null, $totalWs)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "2022-Q4"

# --- 2. Overwrite the fund holdings data with the 2022-Q4 figures ---
$rows = @(
    @("161222", "国投瑞银瑞利灵活配置混合（LOF）A",   "24.37", "75.57", "2.49", "0.6068", 6),
    @("010338", "国投瑞银远见成长混合A",               "10.16", "86.73", "2.56", "0.2601", 7),
    @("121010", "国投瑞银瑞源灵活配置混合A",           "9.28",  "72.64", "2.37", "0.2199", 5),
    @("015652", "国投瑞银瑞利灵活配置混合（LOF）C",   "8.80",  "75.57", "2.49", "0.2191", 6),
    @("010339", "国投瑞银远见成长混合C",               "2.06",  "86.73", "2.56", "0.0527", 7),
    @("015572", "国投瑞银瑞源灵活配置混合C",           "1.90",  "72.64", "2.37", "0.0450", 5),
    @("014541", "华安新能源主题混合A",                 "1.09",  "90.05", "3.08", "0.0336", 9),
    @("620001", "金元顺安宝石动力混合",                 "0.46",  "56.89", "6.47", "0.0298", 3),
    @("015564", "大成弘远回报一年持有混合A",           "2.54",  "27.63", "0.81", "0.0206", 9),
    @("014542", "华安新能源主题混合C",                 "0.10",  "90.05", "3.08", "0.0031", 9),
    @("015565", "大成弘远回报一年持有混合C",           "0.09",  "27.63", "0.81", "0.0007", 9)
)

# fund code / name / size / position / ratio / value columns are stored as
# text in this workbook (they keep trailing & leading zeros), rank stays numeric
$newWs.Range("B2:G12").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newWs.Cells.Item($r, 2).Value = $row[0]
    $newWs.Cells.Item($r, 3).Value = $row[1]
    $newWs.Cells.Item($r, 4).Value = $row[2]
    $newWs.Cells.Item($r, 5).Value = $row[3]
    $newWs.Cells.Item($r, 6).Value = $row[4]
    $newWs.Cells.Item($r, 7).Value = $row[5]
    $newWs.Cells.Item($r, 8).Value = $row[6]
}

# --- 3. Refresh the "总计" summary sheet (2022-Q4 row + shifted history) ---
$summary = @(
    @(0, "2022-Q4", 11, 1.49),
    @(1, "2022-Q3", 11, 1.55),
    @(2, "2022-Q2", 10, 2.19),
    @(3, "2022-Q1", 19, 3.77),
    @(4, "2021-Q4", 17, 4.26),
    @(5, "2021-Q3", 7,  0.86),
    @(6, "2021-Q2", 5,  0.31),
    @(7, "2021-Q1", 8,  0.68),
    @(8, "2020-Q4", 6,  0.52)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $r = $i + 2
    $row = $summary[$i]
    $totalWs.Cells.Item($r, 1).Value = $row[0]
    $totalWs.Cells.Item($r, 2).Value = $row[1]
    $totalWs.Cells.Item($r, 3).Value = $row[2]
    $totalWs.Cells.Item($r, 4).Value = $row[3]
}
